# Updates after KYCC and MIT classes.
#
# 1) Refresh the cached "datetimeFigureOut" date placeholder text
#    (11/27/2018 -> 2/7/2019) on the Slide Master and on every Custom
#    Layout that carries a Date Placeholder shape.
# 2) Bump the doc-number footer textbox on slide 1 from
#    "002-23599 *D" to "002-23599 *E".

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes, $newDate) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $shp = $shapes.Item($i)
        if ($shp.Name -like "Date Placeholder*") {
            $shp.TextFrame.TextRange.Text = $newDate
        }
    }
}

$newDate = "2/7/2019"

# Slide Master
$master = $p.SlideMaster
Update-DatePlaceholder $master.Shapes $newDate

# Every Custom Layout hanging off the master
$layouts = $master.CustomLayouts
for ($li = 1; $li -le $layouts.Count; $li++) {
    $layout = $layouts.Item($li)
    Update-DatePlaceholder $layout.Shapes $newDate
}

# Slide 1: bump the revision tag in the small footer textbox.
$slide = $p.Slides.Item(1)
$shapes = $slide.Shapes
for ($i = 1; $i -le $shapes.Count; $i++) {
    $shp = $shapes.Item($i)
    $txt = $shp.TextFrame.TextRange.Text
    if ($txt -like "002-23599*") {
        $shp.TextFrame.TextRange.Text = "002-23599 *E"
    }
}
